$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.08
$summary.Range("B6").Value = 19
$summary.Range("B9").Value = 26.32

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 19
$status.Range("G4").Value = 26.32

# --- New trade row data (row 20) shared by "All Trades" and "MarketMaking" sheets ---
function Add-TradeRow($ws, $rowIndex) {
    $ws.Cells.Item($rowIndex, 1).Value = 19
    # Leading apostrophe forces text entry so the date-like string "2026-02-17"
    # is not auto-converted into a date serial number by Excel.
    $ws.Cells.Item($rowIndex, 2).Value = "'2026-02-17"
    $ws.Cells.Item($rowIndex, 3).Value = "07:59:36"
    $ws.Cells.Item($rowIndex, 4).Value = "MarketMaking"
    $ws.Cells.Item($rowIndex, 5).Value = "DOWN"
    $ws.Cells.Item($rowIndex, 6).Value = 0.97
    $ws.Cells.Item($rowIndex, 7).Value = 0.97
    $ws.Cells.Item($rowIndex, 8).Value = "CLOSED"
    $ws.Cells.Item($rowIndex, 9).Value = 0
    $ws.Cells.Item($rowIndex, 10).Value = 0
    $ws.Cells.Item($rowIndex, 11).Value = 99.92
    $ws.Cells.Item($rowIndex, 12).Value = 0
    $ws.Cells.Item($rowIndex, 13).Value = 0
    $ws.Cells.Item($rowIndex, 14).Value = 0.6
    $ws.Cells.Item($rowIndex, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($rowIndex, 16).Value = "early_exit"
    $ws.Cells.Item($rowIndex, 17).Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 20

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 20
